# Katalog guncellendi - Pzt 24.11.2025 16:21:53,79
# Adds 5 new "SISME YELEK 5020 NOVA" product rows (KREM, BEJ, KAHVERENGI, SIYAH, TAS)
# at the bottom of the product list (rows 80-84).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the current data
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$startRow = $lastRow + 1

$price = "500 TL"
$category = "Yelek"
$desc = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$stock = "Var"

$r = $startRow

# Row 1 - KREM: name, price, category, image (natural left-to-right order)
$ws.Cells.Item($r, 1).Value = "ŞİŞME YELEK 5020 NOVA KREM"
$ws.Cells.Item($r, 2).Value = $price
$ws.Cells.Item($r, 3).Value = $category
$ws.Cells.Item($r, 4).Value = "5020 NOVA YELEK KREM.jpg"
$ws.Cells.Item($r, 5).Value = $desc
$ws.Cells.Item($r, 6).Value = $stock
$r++

# Row 2 - BEJ: image filename was entered before the product name
$ws.Cells.Item($r, 4).Value = "5020novayelekbej.jpg"
$ws.Cells.Item($r, 1).Value = "ŞİŞME YELEK 5020 NOVA BEJ"
$ws.Cells.Item($r, 2).Value = $price
$ws.Cells.Item($r, 3).Value = $category
$ws.Cells.Item($r, 5).Value = $desc
$ws.Cells.Item($r, 6).Value = $stock
$r++

# Row 3 - KAHVERENGI: natural left-to-right order
$ws.Cells.Item($r, 1).Value = "ŞİŞME YELEK 5020 NOVA KAHVERENGİ"
$ws.Cells.Item($r, 2).Value = $price
$ws.Cells.Item($r, 3).Value = $category
$ws.Cells.Item($r, 4).Value = "5020NOVAYELEKKAHVERENGİ.jpg"
$ws.Cells.Item($r, 5).Value = $desc
$ws.Cells.Item($r, 6).Value = $stock
$r++

# Row 4 - SIYAH: natural left-to-right order
$ws.Cells.Item($r, 1).Value = "ŞİŞME YELEK 5020 NOVA SİYAH"
$ws.Cells.Item($r, 2).Value = $price
$ws.Cells.Item($r, 3).Value = $category
$ws.Cells.Item($r, 4).Value = "5020novayeleksiyah.jpg"
$ws.Cells.Item($r, 5).Value = $desc
$ws.Cells.Item($r, 6).Value = $stock
$r++

# Row 5 - TAS: image filename was entered before the product name
$ws.Cells.Item($r, 4).Value = "5020NOVAYELEKTAŞ.jpg"
$ws.Cells.Item($r, 1).Value = "ŞİŞME YELEK 5020 NOVA TAŞ"
$ws.Cells.Item($r, 2).Value = $price
$ws.Cells.Item($r, 3).Value = $category
$ws.Cells.Item($r, 5).Value = $desc
$ws.Cells.Item($r, 6).Value = $stock

$ws.Range("D82").Select()

$wb.Save()
